# Update the "panel_query_time" timestamps recorded on the "data" sheet.
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

$data.Range("F2").Value = "2021-10-05 14:21:19.422486"
$data.Range("F3").Value = "2021-10-05 14:21:19.422493"
$data.Range("F4").Value = "2021-10-05 14:21:19.422505"
$data.Range("F5").Value = "2021-10-05 14:21:19.422508"
$data.Range("F6").Value = "2021-10-05 14:21:19.422511"

# Add a new "metadata" tab right after "data", carrying the panel query
# metadata that used to be folded into the data sheet.
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Juvenile dermatomyositis"
$meta.Range("C2").Value = 239
$meta.Range("E2").Value = "2020-05-07T14:28:15.158956Z"
$meta.Range("F2").Value = "2021-10-05 14:21:19.418747"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/239/?format=json"

# "data_version" (1.3) must be stored as text, not the number 1.3 -- stage it
# in a scratch cell formatted as text, then paste only the value across so
# the target cell keeps its plain, unformatted style.
$scratch = $meta.Cells.Item(500, 500)
$scratch.NumberFormat = "@"
$scratch.Value = "1.3"
$scratch.Copy()
$meta.Range("D2").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = 0

# Mirror the header style used on the "data" sheet's header row (bold,
# bordered, centered) for the new metadata header row, and the index-column
# style for A2 -- copy the formatting straight from the "data" sheet so the
# same style record is reused.
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
